# Apply crypto price/volume updates per the commit diff.
# Column D values are plain numeric-looking text (e.g. "570.77", "63.778.31")
# that must stay stored as TEXT (inline string), matching the source file.
# A bare $ws.Range(...).Value = "570.77" gets auto-coerced to a Number by
# Excel's COM layer, so we force text via a leading apostrophe (the normal
# Excel "treat as text" convention), then reset .Style so no stray
# quote-prefix formatting lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.046.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "'3.413.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'570.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "'159.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'3.415.58"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").Value = "'0.553"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.86%  "
$ws.Range("D10").Value = "'7.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").Value = "'0.120"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("D12").Value = "'0.425"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.64%  "
$ws.Range("D13").Value = "'4.003.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "'27.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").Value = "'0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.00%  "
$ws.Range("D17").Value = "'64.151.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "'3.462.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "'6.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.35%  "
$ws.Range("D20").Value = "'13.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.66%  "
$ws.Range("D21").Value = "'378.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").Value = "'7.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "'71.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").Value = "'0.517"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.03%  "
$ws.Range("D26").Value = "'0.0000116"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.01%  "
$ws.Range("D27").Value = "'9.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.51%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D30").Value = "'6.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("D31").Value = "'1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.59%  "
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").Value = "'22.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "'7.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("E35").Value = "  -5.95%  "
$ws.Range("D36").Value = "'159.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "'0.852"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.33%  "
$ws.Range("D38").Value = "'1.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.18%  "
$ws.Range("D39").Value = "'2.820.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.93%  "
$ws.Range("D40").Value = "'0.0727"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.20%  "
$ws.Range("D41").Value = "'43.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").Value = "'6.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.64%  "
$ws.Range("D43").Value = "'25.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.33%  "
$ws.Range("D44").Value = "'26.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").Value = "'4.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("D46").Value = "'0.0303"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.04%  "
$ws.Range("D47").Value = "'336.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.59%  "
$ws.Range("D48").Value = "'2.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.16%  "
$ws.Range("D49").Value = "'1.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'6.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.103"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.65%  "
